$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 41986
$ws.Range("B8").Value = 2.1
$ws.Range("C8").Value = "Start styling page header. logo, global menu, site name, site slogan, global header and local header."

$ws.Range("C12").Select() | Out-Null

Write-Output "done"
